$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, pushing the old row 35 (and below) down to row 36.
$ws.Rows.Item(35).Insert()

# Copy the date-cell style (used by column D) from the row above into the new D35 cell.
$ws.Range("D34").Copy()
$ws.Range("D35").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New row 35 gets the values that used to be in row 34 (before this edit), i.e. the
# older "44305" data point, unchanged.
$ws.Range("A35").Value = 9
$ws.Range("B35").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C35").Value = "Metropolitana"
$ws.Range("D35").Value = 44305
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100108
$ws.Range("H35").Value = "Tropicales y subtropicales"
$ws.Range("I35").Value = 100108007
$ws.Range("J35").Value = "Coco"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 40
$ws.Range("N35").Value = 24000
$ws.Range("O35").Value = 24000
$ws.Range("P35").Value = 24000
$ws.Range("Q35").Value = '$/malla 20 unidades'
$ws.Range("R35").Value = "Perú"
$ws.Range("S35").Value = 1200
$ws.Range("T35").Value = 20

# Row 34 is updated with the new, most recent data point.
$ws.Range("D34").Value = 44522
$ws.Range("M34").Value = 25
$ws.Range("N34").Value = 30000
$ws.Range("O34").Value = 30000
$ws.Range("P34").Value = 30000
$ws.Range("S34").Value = 1500

# Update the used range / dimension reflects the extra row automatically.
